$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F, re-using the same header style (bold, bordered,
# centered) as the other header cells in row 1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "Trening"

# New rows of data inserted "above" the previous data (Duza Gra session),
# written directly into rows 2-7 (overwriting what used to be there).
$newRows = @(
    @(45684.59316157408, 649.1, 13.19, 1.967146004949297, "10-15", "Duża Gra"),
    @(45684.59340578704, 670.2, 13.85, 1.988160729408263, "10-15", "Duża Gra"),
    @(45684.59368240741, 694.1, 10.32, 2.056255766323634, "10-15", "Duża Gra"),
    @(45684.59138148148, 495.3, 7.14, 1.939501081194197, "5-10", "Duża Gra"),
    @(45684.59367893518, 693.8, 8.58, 2.277262892041886, "5-10", "Duża Gra"),
    @(45684.59392430555, 715, 7.43, 2.21865485395704, "5-10", "Duża Gra")
)

# The previous data (Mala Gra session), now pushed down to rows 8-13.
$oldRows = @(
    @(45684.59473680556, 785.2, 12.28, 3.682137421199255, "10-15", "Mała Gra"),
    @(45684.59501226852, 809, 11.98, 3.652564287185669, "10-15", "Mała Gra"),
    @(45684.59510023148, 816.6, 12.08, 3.816800253731864, "10-15", "Mała Gra"),
    @(45684.594525, 766.9, 9.140000000000001, 3.410823447363717, "5-10", "Mała Gra"),
    @(45684.59473449074, 785, 9.57, 3.468593393053328, "5-10", "Mała Gra"),
    @(45684.59509791667, 816.4, 9.140000000000001, 3.569871187210084, "5-10", "Mała Gra")
)

$r = 2
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

foreach ($row in $oldRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Apply the datetime number format to the whole Timestamp column (A2:A13).
# Applied twice (lower-case then upper-case) on a single cell first, to
# mirror the author's edit, which registers both custom formats (164 then
# 165) but keeps only the final one (165) referenced by the cell style.
# Then copy that same resulting style to the rest of the column so every
# cell shares one style entry instead of generating a new one per cell.
$ws.Cells.Item(2, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
for ($i = 3; $i -le 13; $i++) {
    $ws.Cells.Item($i, 1).NumberFormat = $ws.Cells.Item(2, 1).NumberFormat
}
